$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Date (column B) values for rows 2-33 with new execution timestamps
$ws.Range("B2").Value = "Mon Sep 11 13:56:22 EDT 2023"
$ws.Range("B3").Value = "Mon Sep 11 13:56:35 EDT 2023"
$ws.Range("B4").Value = "Mon Sep 11 13:56:48 EDT 2023"
$ws.Range("B5").Value = "Mon Sep 11 13:57:02 EDT 2023"
$ws.Range("B6").Value = "Mon Sep 11 13:57:15 EDT 2023"
$ws.Range("B7").Value = "Mon Sep 11 13:57:27 EDT 2023"
$ws.Range("B8").Value = "Mon Sep 11 13:57:40 EDT 2023"
$ws.Range("B9").Value = "Mon Sep 11 13:57:53 EDT 2023"
$ws.Range("B10").Value = "Mon Sep 11 13:58:06 EDT 2023"
$ws.Range("B11").Value = "Mon Sep 11 13:58:19 EDT 2023"
$ws.Range("B12").Value = "Mon Sep 11 13:58:32 EDT 2023"
$ws.Range("B13").Value = "Mon Sep 11 13:58:45 EDT 2023"
$ws.Range("B14").Value = "Mon Sep 11 13:58:58 EDT 2023"
$ws.Range("B15").Value = "Mon Sep 11 13:59:11 EDT 2023"
$ws.Range("B16").Value = "Mon Sep 11 13:59:24 EDT 2023"
$ws.Range("B17").Value = "Mon Sep 11 13:59:37 EDT 2023"
$ws.Range("B18").Value = "Mon Sep 11 13:59:50 EDT 2023"
$ws.Range("B19").Value = "Mon Sep 11 14:00:03 EDT 2023"
$ws.Range("B20").Value = "Mon Sep 11 14:00:17 EDT 2023"
$ws.Range("B21").Value = "Mon Sep 11 14:00:29 EDT 2023"
$ws.Range("B22").Value = "Mon Sep 11 14:00:42 EDT 2023"
$ws.Range("B23").Value = "Mon Sep 11 14:00:56 EDT 2023"
$ws.Range("B24").Value = "Mon Sep 11 14:01:08 EDT 2023"
$ws.Range("B25").Value = "Mon Sep 11 14:01:21 EDT 2023"
$ws.Range("B26").Value = "Mon Sep 11 14:01:34 EDT 2023"
$ws.Range("B27").Value = "Mon Sep 11 14:01:47 EDT 2023"
$ws.Range("B28").Value = "Mon Sep 11 14:02:01 EDT 2023"
$ws.Range("B29").Value = "Mon Sep 11 14:02:14 EDT 2023"
$ws.Range("B30").Value = "Mon Sep 11 14:02:27 EDT 2023"
$ws.Range("B31").Value = "Mon Sep 11 14:02:40 EDT 2023"
$ws.Range("B32").Value = "Mon Sep 11 14:02:53 EDT 2023"
$ws.Range("B33").Value = "Mon Sep 11 14:03:27 EDT 2023"

# Row 32 Result flipped from Pass to Fail on this re-run
$ws.Range("A32").Value = "Fail"

